# Briefing/metrics update for the week of 8/6 (column G), applied to Sheet1.
# - "Pages"  (row 4)  jumped to 158 (was tracking flat at 47).
# - "SLOC (Physical)"/"SLOC (Logical)" (rows 8/9) carried forward unchanged
#   from last week's column F into the new week's column G.
# - "# of Source Files" (row 11) ticked up to 112.
# - "Control Elements Prototyped" (row 15) ticked up to 2.
# Finally, leave the sheet scrolled/selected where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G4").Value = 158
$ws.Range("G8").Value = 9574
$ws.Range("G9").Value = 7025
$ws.Range("G11").Value = 112
$ws.Range("G15").Value = 2

# Restore the cursor/scroll position the author ended up at.
$ws.Range("G23").Select() | Out-Null
